$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 33/34, pushing the existing rows 33.. down to 35..
$ws.Rows("33:34").Insert()

# Row 33 - new "Primera" record for date 44708
$ws.Cells.Item(33, 1).Value = 1
$ws.Cells.Item(33, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(33, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(33, 4).Value = 44708
$ws.Cells.Item(33, 5).Value = 15
$ws.Cells.Item(33, 6).Value = 100112036
$ws.Cells.Item(33, 7).Value = "Caigua"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 120
$ws.Cells.Item(33, 11).Value = 9000
$ws.Cells.Item(33, 12).Value = 10000
$ws.Cells.Item(33, 13).Value = 9500
$ws.Cells.Item(33, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(33, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(33, 16).Value = 475
$ws.Cells.Item(33, 17).Value = 20
$ws.Cells.Item(33, 18).Value = "Hortaliza"

# Row 34 - new "Segunda" record for date 44708
$ws.Cells.Item(34, 1).Value = 1
$ws.Cells.Item(34, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(34, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(34, 4).Value = 44708
$ws.Cells.Item(34, 5).Value = 15
$ws.Cells.Item(34, 6).Value = 100112036
$ws.Cells.Item(34, 7).Value = "Caigua"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Segunda"
$ws.Cells.Item(34, 10).Value = 120
$ws.Cells.Item(34, 11).Value = 7000
$ws.Cells.Item(34, 12).Value = 8000
$ws.Cells.Item(34, 13).Value = 7500
$ws.Cells.Item(34, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(34, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(34, 16).Value = 375
$ws.Cells.Item(34, 17).Value = 20
$ws.Cells.Item(34, 18).Value = "Hortaliza"
